$d = $word.ActiveDocument

$replacements = @(
    @("87÷5=17, 2", "35÷4=8, 3"),
    @("61÷2=30, 1", "64÷3=21, 1"),
    @("22÷6=3, 4", "76÷3=25, 1"),
    @("54÷6=9, 0", "15÷8=1, 7"),
    @("69÷5=13, 4", "42÷8=5, 2"),
    @("15÷7=2, 1", "39÷9=4, 3"),
    @("76÷5=15, 1", "46÷9=5, 1"),
    @("65÷9=7, 2", "51÷2=25, 1"),
    @("68÷4=17, 0", "47÷8=5, 7"),
    @("41÷6=6, 5", "58÷7=8, 2"),
    @("38÷7=5, 3", "69÷9=7, 6"),
    @("96÷8=12, 0", "34÷2=17, 0"),
    @("48÷5=9, 3", "10÷6=1, 4"),
    @("13÷7=1, 6", "30÷9=3, 3"),
    @("16÷6=2, 4", "22÷6=3, 4"),
    @("44÷6=7, 2", "48÷7=6, 6"),
    @("37÷7=5, 2", "94÷6=15, 4"),
    @("16÷5=3, 1", "99÷7=14, 1"),
    @("37÷6=6, 1", "80÷9=8, 8"),
    @("10÷2=5, 0", "60÷5=12, 0"),
    @("16÷8=2, 0", "88÷5=17, 3"),
    @("33÷3=11, 0", "77÷4=19, 1"),
    @("65÷5=13, 0", "45÷6=7, 3"),
    @("48÷8=6, 0", "44÷7=6, 2"),
    @("23÷4=5, 3", "20÷9=2, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Write-Output "done"
